$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Andrew Case's logged time (+5.5 hours: 79h 15m -> 84h 45m)
$ws.Range("B4").Value = "84h 45m"

# Move the active selection to B4 to match the saved workbook view
$ws.Range("B4").Select()
